$wb = $excel.ActiveWorkbook

# Sheet "展览" - row 4 (event "南宁·万圣漫控嘉年华10") and row 5 (event "南宁·梦中礼Lolita茶会")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 670
$ws1.Range("F5").Value = 57

# Sheet "全部类型" - row 5 (event "南宁·万圣漫控嘉年华10") and row 6 (event "南宁·梦中礼Lolita茶会")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 670
$ws4.Range("F6").Value = 57
